$d = $word.ActiveDocument

$replacements = @(
    @("2024-08-29 Thursday", "2024-08-30 Friday"),
    @("170÷9=18, 8", "488÷4=122, 0"),
    @("575÷2=287, 1", "393÷3=131, 0"),
    @("399÷6=66, 3", "459÷4=114, 3"),
    @("413÷9=45, 8", "847÷8=105, 7"),
    @("258÷5=51, 3", "394÷4=98, 2"),
    @("154÷5=30, 4", "288÷8=36, 0"),
    @("425÷6=70, 5", "638÷6=106, 2"),
    @("505÷8=63, 1", "207÷9=23, 0"),
    @("699÷6=116, 3", "958÷2=479, 0"),
    @("137÷4=34, 1", "720÷8=90, 0"),
    @("499÷3=166, 1", "846÷8=105, 6"),
    @("685÷3=228, 1", "204÷5=40, 4"),
    @("844÷3=281, 1", "731÷5=146, 1"),
    @("581÷4=145, 1", "700÷8=87, 4"),
    @("978÷4=244, 2", "152÷3=50, 2"),
    @("785÷2=392, 1", "758÷5=151, 3"),
    @("603÷3=201, 0", "526÷2=263, 0"),
    @("270÷4=67, 2", "686÷6=114, 2"),
    @("492÷9=54, 6", "301÷5=60, 1"),
    @("508÷6=84, 4", "316÷9=35, 1"),
    @("855÷7=122, 1", "439÷8=54, 7"),
    @("946÷6=157, 4", "923÷9=102, 5"),
    @("420÷5=84, 0", "111÷9=12, 3"),
    @("120÷2=60, 0", "404÷4=101, 0"),
    @("146÷5=29, 1", "903÷6=150, 3")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
